# Support for each database platform.
# Adds a second sample data row to the "all_type" sheet:
#   - moves the date-formatted empty cell from C2 down to C3
#   - keeps F2's existing value (1) and adds F3 = 2
#   - adds a text-formatted time value in J2 ("time_type" column)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move C2's style (date number format) down to C3, then remove C2 entirely.
$ws.Range("C2").Copy()
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("C2").Clear()

# New second data row: F3 = 2
$ws.Range("F3").Value = 2

# New J2 cell: text-formatted time-of-day string for the time_type column
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "12:13:14.987654321"
